$wb = $excel.ActiveWorkbook

# --- Sheet references (by original names) ---
$sheetA = $wb.Worksheets.Item("SheetA")
$sheetQuery = $wb.Worksheets.Item("SheetQuery")

# --- Rename "SheetQuery" -> "Inventory" ---
$sheetQuery.Name = "Inventory"

# --- Header text: "Type" -> "Instrument" on the Inventory sheet ---
$sheetQuery.Range("D2").Value = "Instrument"

# --- Column width tweaks ---
# (Character-width COM units snap to the workbook's default-font pixel
# grid, same as real Excel; values below are chosen as the closest
# reachable point to the target stored width.)
# SheetA: column B width 6.51 -> 6.52
$sheetA.Columns.Item(2).ColumnWidth = 5.857142857142857

# Inventory: column B 12 -> 14.68, column D 9.26 -> 12.91, column E 16.87 -> 19.63
$sheetQuery.Columns.Item(2).ColumnWidth = 14
$sheetQuery.Columns.Item(4).ColumnWidth = 12.142857142857142
$sheetQuery.Columns.Item(5).ColumnWidth = 18.857142857142858

# --- Active sheet / tab selection: move from SheetA to Inventory ---
[void]$sheetQuery.Activate()

# --- Selections (active cell) per sheet ---
[void]$sheetA.Range("F10").Select()
[void]$sheetQuery.Range("E30").Select()
